$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.648.86"
$ws.Range("E2").Value = "  +4.83%  "
$ws.Range("D3").Value = "2.732.17"
$ws.Range("E3").Value = "  +3.44%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "578.53"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").Value = "157.56"
$ws.Range("E6").Value = "  +9.53%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "0.612"
$ws.Range("E8").Value = "  +2.09%  "
$ws.Range("D9").Value = "2.757.42"
$ws.Range("E9").Value = "  +3.74%  "
$ws.Range("D10").Value = "6.71"
$ws.Range("E10").Value = "  +2.16%  "
$ws.Range("D11").Value = "0.112"
$ws.Range("E11").Value = "  +5.61%  "
$ws.Range("E12").Value = "  +4.00%  "
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("D14").Value = "3.249.61"
$ws.Range("E14").Value = "  +4.48%  "
$ws.Range("D15").Value = "27.49"
$ws.Range("E15").Value = "  +5.24%  "
$ws.Range("D16").Value = "63.670.09"
$ws.Range("E16").Value = "  +4.91%  "
$ws.Range("D17").Value = "0.0000154"
$ws.Range("E17").Value = "  +8.07%  "
$ws.Range("D18").Value = "2.763.36"
$ws.Range("E18").Value = "  +4.13%  "
$ws.Range("D19").Value = "12.08"
$ws.Range("E19").Value = "  +4.70%  "
$ws.Range("D20").Value = "4.93"
$ws.Range("E20").Value = "  +4.51%  "
$ws.Range("D21").Value = "361.88"
$ws.Range("E21").Value = "  +3.45%  "
$ws.Range("D22").Value = "6.98"
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("D25").Value = "67.01"
$ws.Range("E25").Value = "  +4.89%  "
$ws.Range("D26").Value = "0.170"
$ws.Range("E26").Value = "  +5.61%  "
$ws.Range("D27").Value = "8.58"
$ws.Range("E27").Value = "  +4.85%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("D29").Value = "0.0₃0911"
$ws.Range("E29").Value = "  +13.46%  "
$ws.Range("D30").Value = "2.03"
$ws.Range("E30").Value = "  +1.62%  "
$ws.Range("D31").Value = "7.19"
$ws.Range("E31").Value = "  +5.87%  "
$ws.Range("D32").Value = "1.26"
$ws.Range("E32").Value = "  +19.96%  "
$ws.Range("D33").Value = "175.70"
$ws.Range("E33").Value = "  +7.51%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("D35").Value = "20.62"
$ws.Range("E35").Value = "  +3.80%  "
$ws.Range("D36").Value = "4.90"
$ws.Range("E36").Value = "  +6.11%  "
$ws.Range("D37").Value = "1.45"
$ws.Range("E37").Value = "  +9.75%  "
$ws.Range("E38").Value = "  +9.71%  "
$ws.Range("E39").Value = "  +11.54%  "
$ws.Range("D40").Value = "4.32"
$ws.Range("E40").Value = "  +6.21%  "
$ws.Range("D41").Value = "338.49"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").Value = "39.36"
$ws.Range("E42").Value = "  +2.69%  "
$ws.Range("D43").Value = "5.79"
$ws.Range("E43").Value = "  +11.56%  "
$ws.Range("D44").Value = "21.87"
$ws.Range("E44").Value = "  +8.03%  "
$ws.Range("D45").Value = "22.28"
$ws.Range("E45").Value = "  +8.53%  "
$ws.Range("D46").Value = "0.0600"
$ws.Range("E46").Value = "  +6.63%  "
$ws.Range("D47").Value = "0.648"
$ws.Range("E47").Value = "  +4.04%  "
$ws.Range("D48").Value = "0.0259"
$ws.Range("E48").Value = "  +4.44%  "
$ws.Range("D49").Value = "137.51"
$ws.Range("E49").Value = "  +3.48%  "
$ws.Range("E50").Value = "  +2.35%  "
$ws.Range("D51").Value = "0.995"
$ws.Range("E51").Value = "  -0.32%  "
